$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.816.10"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +5.29%  '
$ws.Range("D3").Value = "'3.107.48"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.95%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'585.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.67%  '
$ws.Range("D6").Value = "'144.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.83%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'3.099.59"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.01%  '
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("E10").Value = '  +11.10%  '
$ws.Range("D11").Value = "'5.70"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +7.68%  '
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("E13").Value = '  +5.32%  '
$ws.Range("D14").Value = "'35.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.99%  '
$ws.Range("D16").Value = "'3.620.93"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.96%  '
$ws.Range("D17").Value = "'7.19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").Value = "'3.103.42"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.91%  '
$ws.Range("D19").Value = "'62.753.10"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.14%  '
$ws.Range("D20").Value = "'462.95"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +6.08%  '
$ws.Range("D21").Value = "'14.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("E23").Value = '  +5.56%  '
$ws.Range("D24").Value = "'13.35"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = "'2.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.69%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +5.46%  '
$ws.Range("D31").Value = "'6.82"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +8.48%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'26.91"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.25%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.110"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.92%  '
$ws.Range("E34").Value = '  +4.82%  '
$ws.Range("E35").Value = '  +11.43%  '
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("D38").Value = "'3.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +12.53%  '
$ws.Range("D39").Value = "'50.92"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.62%  '
$ws.Range("D40").Value = "'8.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.76%  '
$ws.Range("D41").Value = "'427.62"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.88%  '
$ws.Range("E42").Value = '  +4.45%  '
$ws.Range("E43").Value = '  +3.53%  '
$ws.Range("E44").Value = '  +9.55%  '
$ws.Range("E45").Value = '  +3.00%  '
$ws.Range("E46").Value = '  +6.90%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = "'34.95"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.11%  '
$ws.Range("D49").Value = "'123.43"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("E51").Value = '  +4.58%  '
